$wb = $excel.ActiveWorkbook

# --- Update "Version History" sheet: add new version row (v1.3) ---
$wsHistory = $wb.Worksheets.Item("Version History")
$wsHistory.Range("C5").Value = "Set Owner Status to closed"
$wsHistory.Range("A5").Value = "v1.3"
$wsHistory.Range("B5").Value = "Hala Eldaly"
$wsHistory.Range("D5").Value = 45768

# --- Update "Review sheet": set Owner Status to Closed for rows 15-20 ---
$wsReview = $wb.Worksheets.Item("Review sheet")
$wsReview.Range("I15:I20").Value = "Closed"

# --- Restore cursor/selection positions to match the author's final view state ---
$wsReview.Activate()
$wsReview.Range("I1").Select()
$excel.ActiveWindow.ScrollRow = 19

$wsHistory.Activate()
$wsHistory.Range("B14").Select()
